$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("D2").Value = 50
$ws1.Range("H2").Value = 5.98
$ws1.Range("L2").Value = 0.88

$ws1.Range("D3").Value = 55
$ws1.Range("H3").Value = 4.51
$ws1.Range("L3").Value = 0.99

$ws1.Range("H4").Value = 4.58
$ws1.Range("L4").Value = 1.05

$ws1.Range("H5").Value = 3.28
$ws1.Range("L5").Value = 0.85

$ws1.Range("H6").Value = 2.35
$ws1.Range("L6").Value = 1.02

$ws1.Range("H7").Value = 1.35

$ws1.Range("H8").Value = 0.35
$ws1.Range("L8").Value = 0.83

$ws1.Range("L10").Value = 1.05

$ws1.Range("L11").Value = 0.98

$ws1.Range("D12").Value = 42
$ws1.Range("L12").Value = 1.15

$ws1.Range("D13").Value = 39
$ws1.Range("L13").Value = 0.98

$ws1.Range("L14").Value = 1.1

$ws1.Range("D15").Value = 33
$ws1.Range("L15").Value = 1.13

$ws1.Range("D16").Value = 33
$ws1.Range("L16").Value = 0.96

$ws1.Range("D17").Value = 29
$ws1.Range("L17").Value = 0.99

# --- Sheet: Summary ---
# These "numeric" values are stored as text (inline/shared strings) in the
# original workbook, so prefix with an apostrophe to force Excel to keep
# them as text rather than coercing to a number.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'682"
$ws2.Range("B10").Value = "'378"
$ws2.Range("B11").Value = "'196"
$ws2.Range("B12").Value = "'56"
$ws2.Range("B14").Value = "'29"
